$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 23:46"

$ws.Range("B7").Value = 23940
$ws.Range("C7").Value = 4557
$ws.Range("E7").Value = 23468
$ws.Range("G7").Value = 45
$ws.Range("H7").Value = 301

$ws.Range("B8").Value = 22364
$ws.Range("C8").Value = 2516
$ws.Range("E8").Value = 22071

$ws.Range("B12").Value = 6863
$ws.Range("C12").Value = 1248
$ws.Range("D12").Value = 131

$ws.Range("F74").Value = 3

$ws.Range("A94").Value = "Guadalupe"
$ws.Range("C94").Value = 5
$ws.Range("D94").Value = 0
$ws.Range("E94").Value = 55
$ws.Range("F94").Value = 4
$ws.Range("H94").Value = 1

$ws.Range("A95").Value = "Senegal"
$ws.Range("B95").Value = 56
$ws.Range("C95").Value = 9
$ws.Range("D95").Value = 5
$ws.Range("E95").Value = 51

$ws.Range("A96").Value = "Kazajistan"
$ws.Range("B96").Value = 54
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 54

$ws.Range("A97").Value = "Camboya"
$ws.Range("C97").Value = 2
$ws.Range("D97").Value = 2
$ws.Range("E97").Value = 51
$ws.Range("H97").Value = 0

$ws.Range("A98").Value = "Azerbaiyan"
$ws.Range("C98").Value = 9
$ws.Range("D98").Value = 11
$ws.Range("E98").Value = 41
$ws.Range("H98").Value = 1

$ws.Range("A99").Value = "Estado de Palestina"
$ws.Range("B99").Value = 53
$ws.Range("C99").Value = 5
$ws.Range("D99").Value = 17
$ws.Range("E99").Value = 36

$ws.Range("A100").Value = "Nueva Zelanda"
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("E100").Value = 52

$ws.Range("A101").Value = "Oman"
$ws.Range("B101").Value = 52
$ws.Range("C101").Value = 4
$ws.Range("D101").Value = 13
$ws.Range("E101").Value = 39
$ws.Range("F101").Value = 0
$ws.Range("H101").Value = 0

$ws.Range("A104").Value = "Reunion"
$ws.Range("C104").Value = 9
$ws.Range("D104").Value = 0
$ws.Range("E104").Value = 47
$ws.Range("H104").Value = 0

$ws.Range("A105").Value = "Ucrania"
$ws.Range("B105").Value = 47
$ws.Range("C105").Value = 6
$ws.Range("D105").Value = 1
$ws.Range("E105").Value = 43
$ws.Range("H105").Value = 3

$ws.Range("B109").Value = 37
$ws.Range("C109").Value = 5
$ws.Range("E109").Value = 36

$ws.Range("A121").Value = "Guyana"
$ws.Range("C121").Value = 3
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0

$ws.Range("A122").Value = "Paraguay"
$ws.Range("C122").Value = 0
$ws.Range("D122").Value = 0
$ws.Range("F122").Value = 1
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 1

$ws.Range("A123").Value = "Monaco"
$ws.Range("C123").Value = 7
$ws.Range("D123").Value = 1
$ws.Range("E123").Value = 17

$ws.Range("A124").Value = "Macao"
$ws.Range("B124").Value = 18
$ws.Range("C124").Value = 1
$ws.Range("D124").Value = 10
$ws.Range("E124").Value = 8

$ws.Range("A125").Value = "Ruanda"
$ws.Range("C125").Value = 0
$ws.Range("E125").Value = 17
$ws.Range("H125").Value = 0

$ws.Range("A126").Value = "Guatemala"
$ws.Range("B126").Value = 17
$ws.Range("C126").Value = 5
$ws.Range("H126").Value = 1

$ws.Range("A127").Value = "Togo"
$ws.Range("C127").Value = 7

$ws.Range("A128").Value = "Montenegro"
$ws.Range("B128").Value = 16
$ws.Range("C128").Value = 2
$ws.Range("E128").Value = 16

$ws.Range("A130").Value = "Polinesia Francesa"
$ws.Range("B130").Value = 15
$ws.Range("C130").Value = 4
$ws.Range("E130").Value = 15

$ws.Range("A131").Value = "Barbados"
$ws.Range("C131").Value = 8
$ws.Range("E131").Value = 14
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 0

$ws.Range("A132").Value = "Kirguistan"
$ws.Range("C132").Value = 8
$ws.Range("D132").Value = 0
$ws.Range("E132").Value = 14

$ws.Range("A133").Value = "Costa de Marfil"
$ws.Range("B133").Value = 14
$ws.Range("C133").Value = 5
$ws.Range("D133").Value = 1
$ws.Range("E133").Value = 13

$ws.Range("A134").Value = "Mauricio"
$ws.Range("B134").Value = 14
$ws.Range("C134").Value = 2
$ws.Range("E134").Value = 13
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 1

$ws.Range("A135").Value = "Maldivas"
$ws.Range("B135").Value = 13
$ws.Range("D135").Value = 3
$ws.Range("E135").Value = 10

$ws.Range("A136").Value = "Mayotte"
$ws.Range("B136").Value = 11
$ws.Range("C136").Value = 4
$ws.Range("E136").Value = 11

$ws.Range("A137").Value = "Mongolia"
$ws.Range("B137").Value = 10
$ws.Range("C137").Value = 4
$ws.Range("E137").Value = 10

$ws.Range("A138").Value = "Gibraltar"
$ws.Range("B138").Value = 10
$ws.Range("D138").Value = 2
$ws.Range("E138").Value = 8

$ws.Range("A139").Value = "Etiopia"
$ws.Range("B139").Value = 9
$ws.Range("E139").Value = 9

$ws.Range("A140").Value = "Kenia"
$ws.Range("B140").Value = 7
$ws.Range("E140").Value = 7

$ws.Range("A141").Value = "Seychelles"
$ws.Range("B141").Value = 7
$ws.Range("C141").Value = 0
$ws.Range("E141").Value = 7

$ws.Range("A143").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("C143").Value = 3

$ws.Range("A144").Value = "Guinea Ecuatorial"
$ws.Range("B144").Value = 6
$ws.Range("C144").Value = 0
$ws.Range("E144").Value = 6

$ws.Range("A145").Value = "Surinam"
$ws.Range("C145").Value = 1
$ws.Range("D145").Value = 0
$ws.Range("E145").Value = 5

$ws.Range("A146").Value = "San Martin (Parte Francesa)"
$ws.Range("C146").Value = 1
$ws.Range("E146").Value = 5
$ws.Range("H146").Value = 0

$ws.Range("A148").Value = "Aruba"
$ws.Range("B148").Value = 5
$ws.Range("D148").Value = 1

$ws.Range("A149").Value = "Bahamas"
$ws.Range("C149").Value = 0

$ws.Range("A150").Value = "Nueva Caledonia"
$ws.Range("C150").Value = 2

$ws.Range("A152").Value = "El Salvador"

$ws.Range("A153").Value = "Cabo Verde"
$ws.Range("C153").Value = 2

$ws.Range("A154").Value = "Liberia"
$ws.Range("C154").Value = 1

$ws.Range("A155").Value = "Namibia"

$ws.Range("A156").Value = "Zimbabue"
$ws.Range("C156").Value = 2

$ws.Range("A157").Value = "San Bartolome"
$ws.Range("C157").Value = 0

$ws.Range("A158").Value = "Republica de Africa Central"
$ws.Range("C158").Value = 0

$ws.Range("A159").Value = "Congo"

$ws.Range("A160").Value = "Islas Caimanes"

$ws.Range("A161").Value = "Curazao"
